$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.925.66'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '''2.092.68'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''233.14'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D8").Value = '''57.58'
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").Value = '''0.0781'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("D12").Value = '''2.390.19'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '''14.47'
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").Value = '''21.14'
$ws.Range("E14").Value = '  +1.64%  '
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").Value = '''2.094.62'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '''37.854.10'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").Value = '''70.87'
$ws.Range("E20").Value = '  +2.00%  '
$ws.Range("D21").Value = '''0.0₃0822'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '''228.61'
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''2.40'
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").Value = '''170.66'
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").Value = '''0.141'
$ws.Range("E27").Value = '  +11.86%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("E30").Value = '  +1.97%  '
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("E32").Value = '  +3.42%  '
$ws.Range("D33").Value = '''0.0629'
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").Value = '''2.52'
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("E36").Value = '  +3.48%  '
$ws.Range("D37").Value = '''3.41'
$ws.Range("E37").Value = '  +5.66%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = '''5.45'
$ws.Range("E39").Value = '  -4.02%  '
$ws.Range("E40").Value = '  +6.64%  '
$ws.Range("D41").Value = '''2.93'
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = '''97.39'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("D43").Value = '''0.0213'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").Value = '''1.454.62'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  +3.36%  '
$ws.Range("D47").Value = '''15.73'
$ws.Range("E47").Value = '  +4.50%  '
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = '''4.03'
$ws.Range("E48").Value = '  -8.45%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''7.40'
$ws.Range("E49").Value = '  +3.29%  '
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").Value = '''2.286.10'
$ws.Range("E51").Value = '  +0.79%  '
